# Refresh the cryptocurrency price/volume table on Sheet1 (columns B-E, rows 2-51)
# with the latest scraped data. Some coins (rows 29-31) also changed rank, so
# their name/link columns are updated along with price/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D='24.490.09'; E='  -0.63%  ' }
    3 = @{ D='1.694.84'; E='  -0.16%  ' }
    4 = @{ D='1.004'; E='  +0.09%  ' }
    5 = @{ D='316.57'; E='  +0.31%  ' }
    6 = @{ E='  +0.04%  ' }
    7 = @{ D='0.3903'; E='  -0.58%  ' }
    8 = @{ D='0.4068'; E='  +0.81%  ' }
    9 = @{ D='1.484'; E='  -2.15%  ' }
    10 = @{ D='1.004'; E='  +0.03%  ' }
    11 = @{ D='53.14'; E='  +0.49%  ' }
    12 = @{ D='0.08809' }
    13 = @{ D='26.45'; E='  +12.12%  ' }
    14 = @{ D='7.463'; E='  +0.08%  ' }
    15 = @{ D='8.232'; E='  +0.07%  ' }
    16 = @{ D='0.00001361'; E='  +3.32%  ' }
    17 = @{ D='1.693.84'; E='  -0.52%  ' }
    18 = @{ E='  -1.47%  ' }
    19 = @{ D='0.07194'; E='  +2.43%  ' }
    20 = @{ D='20.52'; E='  +4.38%  ' }
    21 = @{ D='7.307'; E='  +3.09%  ' }
    22 = @{ D='1.002'; E='  -0.24%  ' }
    23 = @{ D='14.33'; E='  -1.96%  ' }
    24 = @{ D='24.489.17'; E='  -0.57%  ' }
    25 = @{ D='3.000'; E='  -3.64%  ' }
    26 = @{ D='2.338'; E='  -0.84%  ' }
    27 = @{ D='22.91'; E='  +1.43%  ' }
    28 = @{ D='168.70'; E='  +3.61%  ' }
    29 = @{ B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='5.558'; E='  +7.96%  ' }
    30 = @{ B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='144.58'; E='  +6.58%  ' }
    31 = @{ B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='8.420'; E='  -2.64%  ' }
    32 = @{ D='1.882.26'; E='  -0.47%  ' }
    33 = @{ D='2.193'; E='  +11.53%  ' }
    34 = @{ D='0.08763'; E='  -2.10%  ' }
    35 = @{ D='1.050'; E='  -1.45%  ' }
    36 = @{ D='7.229'; E='  -4.50%  ' }
    37 = @{ E='  +7.96%  ' }
    38 = @{ D='0.2805'; E='  +1.88%  ' }
    39 = @{ D='10.91'; E='  -1.50%  ' }
    40 = @{ D='0.09175'; E='  +0.62%  ' }
    41 = @{ D='14.20'; E='  -1.52%  ' }
    42 = @{ D='0.7938'; E='  +3.88%  ' }
    43 = @{ E='  +2.00%  ' }
    44 = @{ D='17.47'; E='  +10.19%  ' }
    45 = @{ D='2.669'; E='  +4.66%  ' }
    46 = @{ D='0.7232'; E='  +1.07%  ' }
    47 = @{ D='4.264'; E='  +1.33%  ' }
    48 = @{ D='1.397'; E='  +3.82%  ' }
    49 = @{ D='1.002'; E='  +0.01%  ' }
    50 = @{ D='140.16'; E='  +0.31%  ' }
    51 = @{ D='0.08203'; E='  +2.88%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        $value = $rowData[$col]

        if ($col -eq 'D') {
            # Price column: many values look numeric (e.g. "1.004", "0.08809")
            # and would otherwise be auto-converted to numbers by Excel, losing
            # the original text formatting. Force the cell to Text before
            # assigning, then clear the number-format override so the cell
            # keeps its original (unstyled) appearance.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
